$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to force literal-text values for percentage strings
# (direct .Value assignment of e.g. "70%" gets auto-parsed by Excel into
# the number 0.7 with a percent format; round-tripping the text through a
# Text-formatted helper cell + PasteSpecial(xlPasteValues) keeps it literal).
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$ws.Range("E2").Value = '2026-02-21 21:18:37'
$ws.Range("E3").Value = '2026-02-21 21:18:40'
$ws.Range("E4").Value = '2026-02-21 21:18:43'
$ws.Range("O4").Value = '9.4 °C'
$ws.Range("E5").Value = '2026-02-21 21:18:45'
$ws.Range("E6").Value = '2026-02-21 21:18:48'
$helper.Value = '70%'
$helper.Copy()
$ws.Range("H6").PasteSpecial(-4163)
$ws.Range("E7").Value = '2026-02-21 21:18:51'
$helper.Value = '56%'
$helper.Copy()
$ws.Range("H7").PasteSpecial(-4163)
$ws.Range("E8").Value = '2026-02-21 21:18:53'
$helper.Value = '62%'
$helper.Copy()
$ws.Range("H8").PasteSpecial(-4163)
$ws.Range("E9").Value = '2026-02-21 21:18:56'
$helper.Value = '56%'
$helper.Copy()
$ws.Range("H9").PasteSpecial(-4163)
$ws.Range("N9").Value = '6.8 °C 20:50 TU'
$ws.Range("O9").Value = '13.3 °C'
$ws.Range("E10").Value = '2026-02-21 21:18:59'
$helper.Value = '79%'
$helper.Copy()
$ws.Range("H10").PasteSpecial(-4163)
$ws.Range("O10").Value = '8.5 °C'
$ws.Range("E11").Value = '2026-02-21 21:19:00'
$ws.Range("O11").Value = '8.9 °C'
$ws.Range("E12").Value = '2026-02-21 21:19:01'
$helper.Value = '62%'
$helper.Copy()
$ws.Range("H12").PasteSpecial(-4163)
$ws.Range("O12").Value = '12.7 °C'
$ws.Range("E13").Value = '2026-02-21 21:19:02'
$ws.Range("J13").Value = '1031.8 hPa'
$ws.Range("E14").Value = '2026-02-21 21:19:03'
$helper.Value = '70%'
$helper.Copy()
$ws.Range("H14").PasteSpecial(-4163)
$ws.Range("E15").Value = '2026-02-21 21:19:05'
$helper.Value = '56%'
$helper.Copy()
$ws.Range("H15").PasteSpecial(-4163)
$ws.Range("O15").Value = '13.1 °C'
$ws.Range("E16").Value = '2026-02-21 21:19:06'
$ws.Range("O16").Value = '2.3 °C'
$ws.Range("E17").Value = '2026-02-21 21:19:07'
$ws.Range("E18").Value = '2026-02-21 21:19:08'
$ws.Range("E19").Value = '2026-02-21 21:19:09'
$ws.Range("E20").Value = '2026-02-21 21:19:10'
$helper.Value = '38%'
$helper.Copy()
$ws.Range("H20").PasteSpecial(-4163)
$ws.Range("E21").Value = '2026-02-21 21:19:12'
$ws.Range("J21").Value = '1030.7 hPa'
$ws.Range("E22").Value = '2026-02-21 21:19:14'
$ws.Range("E23").Value = '2026-02-21 21:19:17'
$ws.Range("E24").Value = '2026-02-21 21:19:20'
$ws.Range("E25").Value = '2026-02-21 21:19:22'
$ws.Range("E26").Value = '2026-02-21 21:19:25'
$ws.Range("J26").Value = '1027.4 hPa'
$ws.Range("E27").Value = '2026-02-21 21:19:28'
$ws.Range("E28").Value = '2026-02-21 21:19:30'
$helper.Value = '73%'
$helper.Copy()
$ws.Range("H28").PasteSpecial(-4163)
$ws.Range("J28").Value = '1029.7 hPa'
$ws.Range("K28").Value = '14.9 MJ/m2'
$ws.Range("E29").Value = '2026-02-21 21:19:33'
$helper.Value = '67%'
$helper.Copy()
$ws.Range("H29").PasteSpecial(-4163)
$ws.Range("N29").Value = '5.8 °C 20:40 TU'
$ws.Range("O29").Value = '11.5 °C'
$ws.Range("E30").Value = '2026-02-21 21:19:35'
$helper.Value = '67%'
$helper.Copy()
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("E31").Value = '2026-02-21 21:19:38'
$ws.Range("O31").Value = '12.3 °C'
$ws.Range("E32").Value = '2026-02-21 21:19:41'
$ws.Range("O32").Value = '5.1 °C'
$ws.Range("E33").Value = '2026-02-21 21:19:43'
$helper.Value = '52%'
$helper.Copy()
$ws.Range("H33").PasteSpecial(-4163)
$ws.Range("E34").Value = '2026-02-21 21:19:46'
$ws.Range("O34").Value = '4.4 °C'
$ws.Range("E35").Value = '2026-02-21 21:19:48'
$ws.Range("J35").Value = '1031.0 hPa'
$ws.Range("O35").Value = '7.6 °C'
$ws.Range("E36").Value = '2026-02-21 21:19:51'
$helper.Value = '58%'
$helper.Copy()
$ws.Range("H36").PasteSpecial(-4163)
$ws.Range("O36").Value = '13.3 °C'
$ws.Range("E37").Value = '2026-02-21 21:19:54'
$helper.Value = '74%'
$helper.Copy()
$ws.Range("H37").PasteSpecial(-4163)
$ws.Range("J37").Value = '1031.5 hPa'
$ws.Range("O37").Value = '5.8 °C'
$ws.Range("E38").Value = '2026-02-21 21:19:57'
$ws.Range("E39").Value = '2026-02-21 21:19:59'
$ws.Range("E40").Value = '2026-02-21 21:20:02'
$helper.Value = '53%'
$helper.Copy()
$ws.Range("H40").PasteSpecial(-4163)
$ws.Range("J40").Value = '1030.6 hPa'
$ws.Range("O40").Value = '8.7 °C'
$ws.Range("E41").Value = '2026-02-21 21:20:05'
$helper.Value = '69%'
$helper.Copy()
$ws.Range("H41").PasteSpecial(-4163)
$ws.Range("O41").Value = '11.3 °C'
$ws.Range("E42").Value = '2026-02-21 21:20:07'
$helper.Value = '75%'
$helper.Copy()
$ws.Range("H42").PasteSpecial(-4163)
$ws.Range("O42").Value = '10.8 °C'
$ws.Range("E43").Value = '2026-02-21 21:20:10'
$ws.Range("E44").Value = '2026-02-21 21:20:12'
$ws.Range("E45").Value = '2026-02-21 21:20:15'
$ws.Range("E46").Value = '2026-02-21 21:20:18'
$helper.Value = '69%'
$helper.Copy()
$ws.Range("H46").PasteSpecial(-4163)
$ws.Range("O46").Value = '9.8 °C'

# Remove the helper column entirely so it leaves no trace in the sheet
$helper.EntireColumn.Delete()
